$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.548.24'
$ws.Range('D2').ClearFormats()

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.83%  '
$ws.Range('E2').ClearFormats()

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.490.55'
$ws.Range('D3').ClearFormats()

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.60%  '
$ws.Range('E3').ClearFormats()

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E4').ClearFormats()

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.49'
$ws.Range('D5').ClearFormats()

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.23%  '
$ws.Range('E5').ClearFormats()

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.25'
$ws.Range('D6').ClearFormats()

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +3.57%  '
$ws.Range('E6').ClearFormats()

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').ClearFormats()

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('E7').ClearFormats()

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.487.73'
$ws.Range('D8').ClearFormats()

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +2.46%  '
$ws.Range('E8').ClearFormats()

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.594'
$ws.Range('D9').ClearFormats()

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +7.31%  '
$ws.Range('E9').ClearFormats()

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.33'
$ws.Range('D10').ClearFormats()

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.48%  '
$ws.Range('E10').ClearFormats()

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +5.89%  '
$ws.Range('E11').ClearFormats()

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.23%  '
$ws.Range('E12').ClearFormats()

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.093.95'
$ws.Range('D13').ClearFormats()

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.47%  '
$ws.Range('E13').ClearFormats()

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.50%  '
$ws.Range('E14').ClearFormats()

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '28.04'
$ws.Range('D15').ClearFormats()

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +4.60%  '
$ws.Range('E15').ClearFormats()

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.568.50'
$ws.Range('D16').ClearFormats()

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.76%  '
$ws.Range('E16').ClearFormats()

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.93%  '
$ws.Range('E17').ClearFormats()

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.500.24'
$ws.Range('D18').ClearFormats()

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.13%  '
$ws.Range('E18').ClearFormats()

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.60%  '
$ws.Range('E19').ClearFormats()

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.95'
$ws.Range('D20').ClearFormats()

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.90%  '
$ws.Range('E20').ClearFormats()

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '388.94'
$ws.Range('D21').ClearFormats()

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.88%  '
$ws.Range('E21').ClearFormats()

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.91'
$ws.Range('D22').ClearFormats()

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.61%  '
$ws.Range('E22').ClearFormats()

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.80'
$ws.Range('D23').ClearFormats()

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.58%  '
$ws.Range('E23').ClearFormats()

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('E24').ClearFormats()

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.91%  '
$ws.Range('E25').ClearFormats()

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +6.17%  '
$ws.Range('E26').ClearFormats()

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.24'
$ws.Range('D27').ClearFormats()

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +8.61%  '
$ws.Range('E27').ClearFormats()

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.10%  '
$ws.Range('E28').ClearFormats()

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.994'
$ws.Range('D29').ClearFormats()

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.50%  '
$ws.Range('E29').ClearFormats()

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.32'
$ws.Range('D30').ClearFormats()

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.88%  '
$ws.Range('E30').ClearFormats()

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.45'
$ws.Range('D31').ClearFormats()

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +4.66%  '
$ws.Range('E31').ClearFormats()

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.70%  '
$ws.Range('E32').ClearFormats()

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.61'
$ws.Range('D33').ClearFormats()

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.69%  '
$ws.Range('E33').ClearFormats()

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.33'
$ws.Range('D34').ClearFormats()

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +4.76%  '
$ws.Range('E34').ClearFormats()

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.60'
$ws.Range('D35').ClearFormats()

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +8.43%  '
$ws.Range('E35').ClearFormats()

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '163.08'
$ws.Range('D36').ClearFormats()

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.03%  '
$ws.Range('E36').ClearFormats()

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.84%  '
$ws.Range('E37').ClearFormats()

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +5.16%  '
$ws.Range('E38').ClearFormats()

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.78'
$ws.Range('D39').ClearFormats()

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +5.92%  '
$ws.Range('E39').ClearFormats()

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0743'
$ws.Range('D40').ClearFormats()

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.95%  '
$ws.Range('E40').ClearFormats()

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +6.25%  '
$ws.Range('E41').ClearFormats()

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '26.34'
$ws.Range('D42').ClearFormats()

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.40%  '
$ws.Range('E42').ClearFormats()

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.783.52'
$ws.Range('D43').ClearFormats()

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.97%  '
$ws.Range('E43').ClearFormats()

$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('B44').ClearFormats()

$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('C44').ClearFormats()

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '26.50'
$ws.Range('D44').ClearFormats()

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.75%  '
$ws.Range('E44').ClearFormats()

$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'OKB'
$ws.Range('B45').ClearFormats()

$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('C45').ClearFormats()

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.79'
$ws.Range('D45').ClearFormats()

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.47%  '
$ws.Range('E45').ClearFormats()

$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('B46').ClearFormats()

$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('C46').ClearFormats()

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.51'
$ws.Range('D46').ClearFormats()

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +4.91%  '
$ws.Range('E46').ClearFormats()

$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'VeChain'
$ws.Range('B47').ClearFormats()

$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C47').ClearFormats()

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0309'
$ws.Range('D47').ClearFormats()

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.08%  '
$ws.Range('E47').ClearFormats()

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '345.66'
$ws.Range('D48').ClearFormats()

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +5.56%  '
$ws.Range('E48').ClearFormats()

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.93%  '
$ws.Range('E49').ClearFormats()

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.65'
$ws.Range('D50').ClearFormats()

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +12.59%  '
$ws.Range('E50').ClearFormats()

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.863'
$ws.Range('D51').ClearFormats()

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +6.18%  '
$ws.Range('E51').ClearFormats()
